$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update pointer file paths (row 20 = dipwell_measurements, row 16 = initial_zeta_pickle, row 7 = sourcesink)
$ws.Range("B20").Value = "data/revised_dipwell_data_from_first_rainfall_record_without_canal_sensors.csv"
$ws.Range("B16").Value = "initial_condition/best_initial_zeta.p"
$ws.Range("B7").Value = "data/sourcesink_dry_year.xlsx"

# Update the active selection to match the saved view state
$ws.Range("B8").Select()
